# Apply data refresh edits to Long27_DataComp sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - United States
$ws.Range("E3").Value = 29921
$ws.Range("F3").Value = 45108

# Row 4 - Euro Area
$ws.Range("E4").Value = 29921
$ws.Range("F4").Value = 45108

# Row 10 - Taiwan
$ws.Range("C10").Value = 499
$ws.Range("F10").Value = 45108

# Row 11 - Canada
$ws.Range("E11").Value = 29891
$ws.Range("F11").Value = 45078

# Row 13 - Switzerland
$ws.Range("C13").Value = 464
$ws.Range("F13").Value = 45108

# Row 17 - Saudi Arabia
$ws.Range("C17").Value = 367
$ws.Range("F17").Value = 45108

# Row 19 - Indonesia
$ws.Range("E19").Value = 29921
$ws.Range("F19").Value = 45108

# Row 21 - Norway
$ws.Range("E21").Value = 29921
$ws.Range("F21").Value = 45108

# Row 24 - Denmark
$ws.Range("C24").Value = 391
$ws.Range("F24").Value = 45108

# Row 28 - Kuwait
$ws.Range("C28").Value = 356
$ws.Range("F28").Value = 45108
